$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.831.20"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.94%  "
$ws.Range("D3").Value = "'1.560.22"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.01%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "'205.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.26%  "
$ws.Range("E6").Value = "  -0.84%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").Value = "'21.60"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.41%  "
$ws.Range("E9").Value = "  +0.18%  "
$ws.Range("E10").Value = "  -0.64%  "
$ws.Range("D11").Value = "'0.0861"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.33%  "
$ws.Range("D12").Value = "'1.781.53"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "'1.578.34"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.04%  "
$ws.Range("E14").Value = "  -0.85%  "
$ws.Range("D15").Value = "'0.513"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.50%  "
$ws.Range("D16").Value = "'26.837.15"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'61.26"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.47%  "
$ws.Range("D18").Value = "'215.13"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.45%  "
$ws.Range("E19").Value = "  +1.72%  "
$ws.Range("E20").Value = "  -0.32%  "
$ws.Range("D22").Value = "'4.13"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.44%  "
$ws.Range("E23").Value = "  -1.93%  "
$ws.Range("D24").Value = "'2.01"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.31%  "
$ws.Range("D25").Value = "'153.25"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.23%  "
$ws.Range("D26").Value = "'6.63"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.13%  "
$ws.Range("D27").Value = "'15.03"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.85%  "
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("E29").Value = "  -0.96%  "
$ws.Range("E30").Value = "  +1.14%  "
$ws.Range("E31").Value = "  -3.51%  "
$ws.Range("D32").Value = "'3.17"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.39%  "
$ws.Range("D33").Value = "'1.380.49"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.41%  "
$ws.Range("D34").Value = "'2.92"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.29%  "
$ws.Range("E35").Value = "  -1.94%  "
$ws.Range("E36").Value = "  -0.76%  "
$ws.Range("D37").Value = "'0.917"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.19%  "
$ws.Range("D38").Value = "'0.0163"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.67%  "
$ws.Range("E39").Value = "  +1.88%  "
$ws.Range("D40").Value = "'0.810"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.10%  "
$ws.Range("E41").Value = "  -0.10%  "
$ws.Range("E42").Value = "  +0.62%  "
$ws.Range("E43").Value = "  +5.02%  "
$ws.Range("D44").Value = "'1.78"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.19%  "
$ws.Range("E45").Value = "  +0.83%  "
$ws.Range("D46").Value = "'63.49"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.41%  "
$ws.Range("D47").Value = "'1.695.12"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("D48").Value = "'86.50"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.40%  "
$ws.Range("D49").Value = "'0.0507"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.07%  "
$ws.Range("D50").Value = "'0.0₇0980"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.60%  "
$ws.Range("E51").Value = "  +1.02%  "
